$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# --- 3. Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy the bold/border/center-top header format from an existing header cell
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# --- 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(44941.99999999999,58,-17.27986140535263,129.5463241526176),
    @(44955.99999999999,60,-18.53725535265832,135.1481297297163),
    @(44969.99999999999,62,-12.51941801081181,133.1746777001956),
    @(44976.99999999999,62,-13.06006631920632,139.4423921446176),
    @(44983.99999999999,63,-15.02995871788602,140.7876882701573),
    @(44990.99999999999,64,-9.525771581496203,141.9402505129319),
    @(44997.99999999999,65,-20.58060758835707,136.5569122175034),
    @(45011.99999999999,67,-9.918800318881321,141.9489912900616),
    @(45025.99999999999,68,0.130333427034417,141.2204776400593),
    @(45039.99999999999,70,-3.445173897806443,138.7538575275177),
    @(45046.99999999999,71,-0.8964615891345376,143.779550791328),
    @(45053.99999999999,72,-3.053883110882889,145.1522163478319),
    @(45060.99999999999,72,-0.4630434193252077,146.6939607897322),
    @(45067.99999999999,73,1.191321690850092,148.1287272311903),
    @(45074.99999999999,74,4.208549685024629,149.6735081289869),
    @(45081.99999999999,75,1.132702505198296,147.9748934388483),
    @(45102.99999999999,77,4.882042156734338,148.7528971654177),
    @(45116.99999999999,79,2.162868622949189,153.9178869861332),
    @(45123.99999999999,80,2.616946573496167,153.0755101902198),
    @(45130.99999999999,81,10.22993384294525,151.8448583241334),
    @(45137.99999999999,82,6.395028541329515,155.30768674424),
    @(45144.99999999999,82,11.02965181845889,156.7438811241125),
    @(45151.99999999999,83,9.833892154798443,161.7323726180264),
    @(45158.99999999999,84,13.24153561553437,157.7078792674662),
    @(45165.99999999999,85,5.708670749312129,161.0206619994737),
    @(45172.99999999999,86,11.1832681130144,160.187551526242),
    @(45179.99999999999,87,8.616638319218762,156.3395773423696),
    @(45186.99999999999,87,15.30972535552784,161.3696660082081),
    @(45193.99999999999,88,19.30155027126166,157.7967851747362),
    @(45200.99999999999,89,14.91058050758159,167.4970339714346)
)

$r = 2
foreach ($row in $data) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Apply the date/time number-format style used in column A on the other sheets
$ws1.Range("A2").Copy()
$ws3.Range("A2:A31").PasteSpecial(-4122)

$ws1.Activate()
$ws1.Range("A1").Select()
